# The deck ships with two embedded themes:
#   theme1.xml -> "Simple Light"  (currently applied to the slide master / all slides)
#   theme2.xml -> "Default"       (currently only referenced by the notes master)
#
# This change recolors the slide master's theme color scheme from the
# "Simple Light" palette to the "Default" palette (teal/green dk2; blue,
# green, orange, yellow, cyan, green accents; blue hyperlink; purple
# followed-hyperlink). That's what changes the look of the case study's
# header/title area.

function ToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

# Theme color slots (1-12): dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$tcs.Colors(1).RGB  = ToRGB "000000"   # dk1
$tcs.Colors(2).RGB  = ToRGB "FFFFFF"   # lt1
$tcs.Colors(3).RGB  = ToRGB "158158"   # dk2
$tcs.Colors(4).RGB  = ToRGB "F3F3F3"   # lt2
$tcs.Colors(5).RGB  = ToRGB "058DC7"   # accent1
$tcs.Colors(6).RGB  = ToRGB "50B432"   # accent2
$tcs.Colors(7).RGB  = ToRGB "ED561B"   # accent3
$tcs.Colors(8).RGB  = ToRGB "EDEF00"   # accent4
$tcs.Colors(9).RGB  = ToRGB "24CBE5"   # accent5
$tcs.Colors(10).RGB = ToRGB "64E572"   # accent6
$tcs.Colors(11).RGB = ToRGB "2200CC"   # hlink
$tcs.Colors(12).RGB = ToRGB "551A8B"   # folHlink
